# New crime data collected - weekly refresh of the 88th Precinct CompStat
# report: header "as of" volume/date runs, and the Week-to-Date / 28-Day /
# Year-to-Date / 2-Year crime statistics table (rows 14-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text runs (rich text: only part of the shared string changes)
# ---------------------------------------------------------------------
# A8 = "Volume 32   Number  17" -> "...18"
$ws.Range("A8").Characters(21, 2).Text = "18"

# C9 = "Report Covering the Week  4/21/2025  Through  4/27/2025"
#   -> "...4/28/2025  Through  5/4/2025"
$ws.Range("C9").Characters(27, 9).Text = "4/28/2025"
$ws.Range("C9").Characters(47, 9).Text = "5/4/2025"

# ---------------------------------------------------------------------
# Crime statistics table - plain numeric value updates (style/type
# unchanged for these cells)
# ---------------------------------------------------------------------
$ws.Range("N14").Value = -90

$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 114.285714285714
$ws.Range("I16").Value = 41
$ws.Range("J16").Value = 31
$ws.Range("K16").Value = 32.258064516129
$ws.Range("L16").Value = 36.666666666666
$ws.Range("M16").Value = -35.9375
$ws.Range("N16").Value = -88.483146067415

$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 40
$ws.Range("I17").Value = 77
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 28.333333333333
$ws.Range("L17").Value = 26.229508196721
$ws.Range("M17").Value = 108.108108108108
$ws.Range("N17").Value = -48.666666666666

$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 7
$ws.Range("H18").Value = -22.222222222222
$ws.Range("I18").Value = 22
$ws.Range("J18").Value = 37
$ws.Range("K18").Value = -40.540540540540
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -45
$ws.Range("N18").Value = -91.164658634538

$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = -4.166666666666
$ws.Range("I19").Value = 105
$ws.Range("J19").Value = 79
$ws.Range("K19").Value = 32.911392405063
$ws.Range("L19").Value = -1.869158878504
$ws.Range("M19").Value = -16
$ws.Range("N19").Value = -38.235294117647

$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 27
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = -46
$ws.Range("L20").Value = 17.391304347826
$ws.Range("M20").Value = -10
$ws.Range("N20").Value = -89.534883720930

$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 7.142857142857
$ws.Range("F21").Value = 74
$ws.Range("G21").Value = 63
$ws.Range("H21").Value = 17.460317460317
$ws.Range("I21").Value = 278
$ws.Range("J21").Value = 259
$ws.Range("K21").Value = 7.335907335907
$ws.Range("L21").Value = 8.171206225680
$ws.Range("M21").Value = -7.641196013289
$ws.Range("N21").Value = -76.948590381426

$ws.Range("L22").Value = 0

$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 51
$ws.Range("K23").Value = 59.375
$ws.Range("L23").Value = 54.545454545454
$ws.Range("M23").Value = 59.375

$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = -40
$ws.Range("F24").Value = 43
$ws.Range("G24").Value = 52
$ws.Range("H24").Value = -17.307692307692
$ws.Range("I24").Value = 198
$ws.Range("J24").Value = 194
$ws.Range("K24").Value = 2.061855670103
$ws.Range("L24").Value = -14.655172413793
$ws.Range("M24").Value = -36.129032258064

$ws.Range("C25").Value = 2
$ws.Range("E25").Value = -50
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = -15.384615384615
$ws.Range("I25").Value = 53
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = -11.666666666666
$ws.Range("L25").Value = -11.666666666666

$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 200
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = -6.25
$ws.Range("I26").Value = 124
$ws.Range("J26").Value = 123
$ws.Range("K26").Value = 0.813008130081
$ws.Range("L26").Value = 11.711711711711
$ws.Range("M26").Value = 27.835051546391

$ws.Range("F28").Value = 3
$ws.Range("L28").Value = 12.5

$ws.Range("I29").Value = 4
$ws.Range("L29").Value = -20
$ws.Range("M29").Value = -42.857142857142
$ws.Range("N29").Value = -88.235294117647

$ws.Range("I30").Value = 4
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -33.333333333333
$ws.Range("N30").Value = -86.206896551724

# ---------------------------------------------------------------------
# Cells that flip between a numeric count and the "no data" text markers
# ("0" / "***.*") used elsewhere in the sheet. Value is set first (with a
# leading apostrophe to force text so short numeric-looking strings like
# "0" are not reinterpreted as numbers), then the number format/style is
# copied in from a cell that already carries the correct style, matching
# how the other text-marker cells in this table are formatted.
# ---------------------------------------------------------------------

# Numeric -> "0" (shared text), right-aligned General style
$ws.Range("C22").Value = "'0"
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("D23").Value = "'0"
$ws.Range("D22").Copy()
$ws.Range("D23").PasteSpecial(-4122)

$ws.Range("C28").Value = "'0"
$ws.Range("D22").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("G28").Value = "'0"
$ws.Range("D22").Copy()
$ws.Range("G28").PasteSpecial(-4122)

# Numeric -> "***.*" (shared text), right-aligned General style
$ws.Range("E23").Value = "'***.*"
$ws.Range("N22").Copy()
$ws.Range("E23").PasteSpecial(-4122)

$ws.Range("H28").Value = "'***.*"
$ws.Range("N22").Copy()
$ws.Range("H28").PasteSpecial(-4122)

# "0" text marker -> numeric 1, #,##0 style
$ws.Range("C29").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C29").PasteSpecial(-4122)

$ws.Range("F29").Value = 1
$ws.Range("C16").Copy()
$ws.Range("F29").PasteSpecial(-4122)

$ws.Range("C30").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$ws.Range("F30").Value = 1
$ws.Range("C16").Copy()
$ws.Range("F30").PasteSpecial(-4122)

$excel.CutCopyMode = $false
